# The presentation currently uses the "Integral" theme (ppt/theme/theme1.xml,
# wired to the slide master that every slide inherits from). The edit swaps
# the deck's applied color theme for the built-in Office default palette -
# the same 12-slot theme color scheme ("Office Theme") that PowerPoint ships
# with out of the box.
#
# The 12 theme colour slots (PowerPoint COM order via ThemeColorScheme):
#   1 dk1=000000  2 lt1=FFFFFF  3 dk2=44546A  4 lt2=E7E6E6
#   5 accent1=5B9BD5  6 accent2=ED7D31  7 accent3=A5A5A5  8 accent4=FFC000
#   9 accent5=4472C4  10 accent6=70AD47  11 hlink=0563C1  12 folHlink=954F72
#
# ThemeColorScheme is a single shared part (the theme backing the slide
# master), so editing it through any one slide updates it for the whole
# deck; we do it through slide 1.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$tcs.Item(1).RGB  = 0x000000   # dk1
$tcs.Item(2).RGB  = 0xFFFFFF   # lt1
$tcs.Item(3).RGB  = 0x6A5444   # dk2      -> 44546A
$tcs.Item(4).RGB  = 0xE6E6E7   # lt2      -> E7E6E6
$tcs.Item(5).RGB  = 0xD59B5B   # accent1  -> 5B9BD5
$tcs.Item(6).RGB  = 0x317DED   # accent2  -> ED7D31
$tcs.Item(7).RGB  = 0xA5A5A5   # accent3  -> A5A5A5
$tcs.Item(8).RGB  = 0x00C0FF   # accent4  -> FFC000
$tcs.Item(9).RGB  = 0xC47244   # accent5  -> 4472C4
$tcs.Item(10).RGB = 0x47AD70   # accent6  -> 70AD47
$tcs.Item(11).RGB = 0xC16305   # hlink    -> 0563C1
$tcs.Item(12).RGB = 0x724F95   # folHlink -> 954F72
